$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.747489333152771
$ws.Range("B1").Value = 1.89954686164856
$ws.Range("C1").Value = 1.816475033760071
$ws.Range("D1").Value = 2.172075986862183
$ws.Range("E1").Value = 3.058873176574707
